# elective sort bug fix
# The "HUS, PPE Applied Year" value in B2 was incorrect; correct it to 2016.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 2016
